$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '79.539.43'
$ws.Range('E2').Value = '  +4.17%  '
$ws.Range('D3').Value = '3.141.92'
$ws.Range('E3').Value = '  +2.15%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '206.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '623.76'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.265'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +23.83%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.582'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.50%  '
$ws.Range('D10').Value = '3.137.98'
$ws.Range('E10').Value = '  +2.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.579'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +30.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000249'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +27.27%  '
$ws.Range('E13').Value = '  +1.41%  '
$ws.Range('D14').Value = '3.717.30'
$ws.Range('E14').Value = '  +2.17%  '
$ws.Range('E15').Value = '  -0.79%  '
$ws.Range('E16').Value = '  +6.80%  '
$ws.Range('D17').Value = '79.295.34'
$ws.Range('E17').Value = '  +4.02%  '
$ws.Range('D18').Value = '3.147.51'
$ws.Range('E18').Value = '  +2.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.11'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.52%  '
$ws.Range('E20').Value = '  +14.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '431.39'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +12.27%  '
$ws.Range('E22').Value = '  -0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.14'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +13.84%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.87'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.53%  '
$ws.Range('D25').Value = '3.305.06'
$ws.Range('E25').Value = '  +2.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '75.28'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.94%  '
$ws.Range('E27').Value = '  +1.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.69'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.00%  '
$ws.Range('E29').Value = '  -0.19%  '
$ws.Range('E30').Value = '  +11.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.996'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.29%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.85'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.42%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '548.94'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +9.66%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.46'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.51%  '
$ws.Range('E35').Value = '  +15.23%  '
$ws.Range('E36').Value = '  +2.40%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '22.74'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +8.94%  '
$ws.Range('E38').Value = '  +18.13%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.997'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.398'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.21%  '
$ws.Range('E41').Value = '  +3.35%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '162.98'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.52'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.97%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '187.62'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.48%  '
$ws.Range('E46').Value = '  +7.63%  '
$ws.Range('E47').Value = '  +8.63%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.776'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.80%  '
$ws.Range('E49').Value = '  +1.31%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '42.38'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.82%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.18'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.17%  '
